# Rows 19-22 on the active sheet ("Artfynd") got their data cyclically
# rotated by the source system: the record that used to live in row 21
# now lives in row 19, the one in row 19 moved to row 20, row 22's moved
# to row 21, and row 20's moved to row 22 (a downward rotation with wrap
# -around). Capture each full row first (A:AY) so the reads all happen
# against the original, un-mutated data before anything is overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row19 = $ws.Range("A19:AY19").Value2
$row20 = $ws.Range("A20:AY20").Value2
$row21 = $ws.Range("A21:AY21").Value2
$row22 = $ws.Range("A22:AY22").Value2

# Columns Y and AA hold dates stored as plain text (e.g. "2023-08-18").
# Forcing the cell to Text format before writing keeps Excel from
# auto-converting the recognizable date string into a date serial value.
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("AA19").NumberFormat = "@"
$ws.Range("Y20").NumberFormat = "@"
$ws.Range("AA20").NumberFormat = "@"
$ws.Range("Y21").NumberFormat = "@"
$ws.Range("AA21").NumberFormat = "@"
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("AA22").NumberFormat = "@"

$ws.Range("A19:AY19").Value = $row21
$ws.Range("A20:AY20").Value = $row19
$ws.Range("A21:AY21").Value = $row22
$ws.Range("A22:AY22").Value = $row20
